# Update the "non-convex generator" experiment values (alpha ~ zero variant).
# Values in these sheets are stored as literal TEXT (shared strings) even
# though many of them look like plain numbers, so we force the target
# ranges to Text format before writing, then clear the number format again
# so the cells fall back onto the workbook's default style (matches the
# original file, which has no custom per-cell number formatting).

$wb = $excel.ActiveWorkbook

function Set-TextValues {
    param($Worksheet, $Addresses, $Values)

    # Union all addresses into one Range so we only touch the stylesheet once.
    $unionRange = $Worksheet.Range($Addresses[0])
    for ($i = 1; $i -lt $Addresses.Count; $i++) {
        $unionRange = $Worksheet.Range($unionRange, $Worksheet.Range($Addresses[$i]))
    }

    $unionRange.NumberFormat = "@"
    for ($i = 0; $i -lt $Addresses.Count; $i++) {
        $Worksheet.Range($Addresses[$i]).Value = $Values[$i]
    }
    $unionRange.ClearFormats()
}

# NOTE: sheet lookup by name is case-insensitive (like real Excel), and this
# workbook has both "Vector_bf" and "Vector_BF" sheets, so we address sheets
# by their 1-based tab position to avoid ambiguity:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

# ---- Restricciones_del_follower ----
$wsFollower = $wb.Worksheets.Item(3)

$followerAddresses = @(
    "A2","B2","D2","E2",
    "A3","B3","D3","E3","F3",
    "A4","B4","D4","E4",
    "A5","B5","D5","E5","F5",
    "A6","B6","D6","E6","F6"
)
$followerValues = @(
    "7.865 - x - 0.5y","-5.865","0.93","0",
    "-4.975 - 0.25x + y","2.9749999999999996","0.41","-0.8","-1.2",
    "-7.865 + x + 0.5y","-0.1349999999999998","0.7","0",
    "-11.559999999999999 + x - 2y","-9.559999999999999","0.36","0","0",
    "-6.17 - y","-6.17","0.79","0","0"
)
Set-TextValues $wsFollower $followerAddresses $followerValues

# ---- Punto_modificado ----
$wsPunto = $wb.Worksheets.Item(4)
Set-TextValues $wsPunto @("A2","B2") @("4.78","6.17")

# ---- Vector_bf ----
$wsVecbf = $wb.Worksheets.Item(5)
Set-TextValues $wsVecbf @("A2") @("2.215")

# ---- Vector_BF ----
$wsVecBF = $wb.Worksheets.Item(6)
Set-TextValues $wsVecBF @("A2","A3") @("-1.2","-0.19999999999999996")
